# Degree of freedom correction for clustered EPA
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row ---
$ws.Range("B1").Value = "Tobs"
$ws.Range("D1").Value = "variant"
$ws.Range("E1").Value = "known"
$ws.Range("F1").Value = "known_cond"
$ws.Range("G1").Value = "naive"
$ws.Range("H1").Value = "naive_cond"
$ws.Range("I1").Value = "split"
$ws.Range("J1").Value = "split_cond"
$ws.Range("K1").Value = "selective"
$ws.Range("L1").Value = "selective_cond"
$ws.Range("M1").Value = "selective_mode2"
$ws.Range("N1").Value = "selective_mode2_cond"

# --- Row 2 (existing row, now gains a "variant" column and new metric columns) ---
$ws.Range("A2").Value = 80
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "overall_holds"
$ws.Range("E2").Value = 0.068
$ws.Range("F2").Value = 0.05
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.999
$ws.Range("I2").Value = 0.068
$ws.Range("J2").Value = 0.109
$ws.Range("K2").Value = 0.054
$ws.Range("L2").Value = 0.054
$ws.Range("M2").Value = 0.067
$ws.Range("N2").Value = 0.05

# --- Row 3 (new) ---
$ws.Range("A3").Value = 80
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "overall_holds"
$ws.Range("E3").Value = 0.049
$ws.Range("F3").Value = 0.042
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0.052
$ws.Range("J3").Value = 0.06
$ws.Range("K3").Value = 0.051
$ws.Range("L3").Value = 0.056
$ws.Range("M3").Value = 0.064
$ws.Range("N3").Value = 0.049

# --- Row 4 (new) ---
$ws.Range("A4").Value = 80
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = "overall_holds"
$ws.Range("E4").Value = 0.061
$ws.Range("F4").Value = 0.056
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.057
$ws.Range("J4").Value = 0.063
$ws.Range("K4").Value = 0.074
$ws.Range("L4").Value = 0.049
$ws.Range("M4").Value = 0.074
$ws.Range("N4").Value = 0.029

# --- Row 5 (new) ---
$ws.Range("A5").Value = 80
$ws.Range("B5").Value = 200
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "overall_holds"
$ws.Range("E5").Value = 0.049
$ws.Range("F5").Value = 0.049
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.034
$ws.Range("J5").Value = 0.047
$ws.Range("K5").Value = 0.046
$ws.Range("L5").Value = 0.053
$ws.Range("M5").Value = 0.049
$ws.Range("N5").Value = 0.041

$wb.Save()
